# ECE 222 LEC 8 edit script
# Strategy: locate each target paragraph (or paragraph span) via unique text
# anchors using Find, then replace the *entire* paragraph range's contents
# with freshly authored OOXML via Range.InsertXML. Replacing a range that
# already spans real text (not a collapsed insertion point) reliably merges
# the new content in place without splitting off stray empty paragraphs.

$d = $word.ActiveDocument

function Get-PkgXml([string]$bodyFragment) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Move the "C - value of specified dest reg" bullet so it follows the
#    "Address C - specify which register to write to" bullet (instead of
#    following "A, B - values of specified source regs"), wrap "dest" in
#    proofErr spellStart/spellEnd, and relocate the _GoBack bookmark onto
#    the end of the "source regs" bullet.
# ---------------------------------------------------------------------

$rStart = $d.Content
$rStart.Find.Execute("Address C – specify", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pStart = $rStart.Paragraphs(1)

$rEnd = $d.Content
$rEnd.Find.Execute("dest reg", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pEnd = $rEnd.Paragraphs(1)

$region1 = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$region1Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:t>Address C – specify which register to write to</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">C – value of specified </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> reg</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:t>Outputs:</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">A, </w:t></w:r><w:r><w:t>B</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">– values of specified </w:t></w:r><w:r><w:t>source regs</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$region1.InsertXML((Get-PkgXml $region1Xml))

# ---------------------------------------------------------------------
# 2) Wrap a collection of single-word bullet paragraphs in
#    proofErr spellStart/spellEnd (Word's "unknown word" spell-check
#    markers). Each of these paragraphs consists of exactly one run.
# ---------------------------------------------------------------------

$simpleWords = @("RF_write", "C_select", "B_select", "ALU_op", "MEM_read", "MEM_write", "IR_enable", "MA_select", "INC_select", "PC_select")

foreach ($word_ in $simpleWords) {
    $fr = $d.Content
    $fr.Find.Execute($word_, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $fr.Paragraphs(1)
    $pRange = $para.Range

    $pPrXml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr>'
    $bodyXml = '<w:p>' + $pPrXml + '<w:proofErr w:type="spellStart"/><w:r><w:t>' + $word_ + '</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

    $pRange.InsertXML((Get-PkgXml $bodyXml))
}

# ---------------------------------------------------------------------
# 3) Split " if PC_enable is on" into " if " + proofErr(PC_enable) + " is on"
#    inside the "Selects between RA & output from adder if PC_enable is on"
#    bullet (first occurrence of PC_enable in the document).
# ---------------------------------------------------------------------

$rSentence = $d.Content
$rSentence.Find.Execute("PC_enable", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$sentencePara = $rSentence.Paragraphs(1)
$sentenceRange = $sentencePara.Range

$sentenceBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Selects between </w:t></w:r>' + `
    '<w:r><w:t>RA</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> &amp; output from adder</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> if </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>PC_enable</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> is on</w:t></w:r>' + `
    '</w:p>'

$sentenceRange.InsertXML((Get-PkgXml $sentenceBody))

# ---------------------------------------------------------------------
# 4) Wrap the standalone "PC_enable" bullet (now the 2nd occurrence, since
#    the sentence above still literally contains the word once) in
#    proofErr spellStart/spellEnd.
# ---------------------------------------------------------------------

$rPc = $d.Content
$rPc.Find.Execute("PC_enable", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rPc.Collapse(0)
$rPc.Find.Execute("PC_enable", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pcPara = $rPc.Paragraphs(1)
$pcRange = $pcPara.Range

$pcBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>PC_enable</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

$pcRange.InsertXML((Get-PkgXml $pcBody))

# ---------------------------------------------------------------------
# 5) Split "Microprogrammed control - i.e. through software" so that
#    "Microprogrammed" is wrapped in proofErr spellStart/spellEnd, AND
#    remove the _GoBack bookmark from the trailing (last) empty paragraph -
#    it was relocated onto the "source regs" bullet in step 1 above, so
#    that last paragraph simply becomes bookmark-free.
#
#    Both paragraphs are replaced together in a single range/InsertXML
#    call: the trailing bookmark paragraph has no real text of its own
#    (only a bookmark), so a range comprising *only* that paragraph sits
#    exactly on a paragraph-mark boundary and InsertXML would split off an
#    extra empty paragraph instead of replacing in place. Folding it into
#    the combined range together with the preceding (non-empty) paragraph
#    avoids that edge case.
# ---------------------------------------------------------------------

$rMicro = $d.Content
$rMicro.Find.Execute("Microprogrammed control", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$microPara = $rMicro.Paragraphs(1)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$microRange = $d.Range($microPara.Range.Start, $lastPara.Range.End)

$microBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Microprogrammed</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> control – i.e. through software</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:p>'

$microRange.InsertXML((Get-PkgXml $microBody))

Write-Output "edit complete"
